# Updated D/D/1 example - Data file modified to reflect proper IAT's
# Data sheet: B column (IAT) now derived from inter-arrival gaps (D_n - D_(n-1))
# and C column formulas re-entered so rows 2-32 share one formula group.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Column C (AT = F + G) re-entered across C2:C32 as one fill, producing a
#     shared formula group anchored at C2. Values are unchanged, only the
#     formula storage becomes "shared".
$ws.Range("C2:C32").Formula = "=F2+G2"

# --- Column B (IAT) switches from hard-coded constants to the gap between
#     consecutive arrival times in column D.
# B2 keeps a literal value.
$ws.Range("B2").Value = 2

# B3 and B4 get their own (non-shared) formulas.
$ws.Range("B3").Formula = "=D3-D2"
$ws.Range("B4").Formula = "=D4-D3"

# B5:B32 are filled as one shared formula group anchored at B5.
$ws.Range("B5:B32").Formula = "=D5-D4"

# --- Update the sheet's active selection from D7 to D5.
$ws.Activate()
$ws.Range("D5").Select()
